$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Start the new row by copying the formatting of the previous data row (A6:C6)
# down into row 7, so the date cell picks up the existing date number format
# instead of Excel registering a brand-new custom number format.
$ws.Range("A6:C6").Copy($ws.Range("A7:C7"))

# Now overwrite with the actual new entry: same date as the last entry,
# 1 additional hour worked, with a new description.
$ws.Range("A7").Value = (Get-Date -Year 2018 -Month 1 -Day 23 -Hour 0 -Minute 0 -Second 0 -Millisecond 0)
$ws.Range("B7").Value = 1
$ws.Range("C7").Value = "more R work"

# Move the active selection down to the next empty row, like Excel does after data entry
$ws.Range("A8").Select()
